$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Foil name relabeling (2 -> 3 suffix), row 2,3,4,5,8 ---
$ws.Range("A2").Value = "Zr3"
$ws.Range("A3").Value = "In3"
$ws.Range("A4").Value = "Ni3"
$ws.Range("A5").Value = "Au3"
$ws.Range("A8").Value = "Al3"

# --- Data corrections ---
# Row 3: Uncertainty E3 0.1 -> 0.01
$ws.Range("E3").Value = 0.01

# Row 5: Density D5 1.001 -> 0.1 (this ripples H5 via existing formula)
$ws.Range("D5").Value = 0.1

# --- New uncertainty column I, rows 2-8 ---
$ws.Range("I2").Formula = "=SQRT((C2/B2)^2+(E2/D2)^2+(G2/F2)^2)*H2"
$ws.Range("I3:I8").Formula = "=SQRT((C3/B3)^2+(E3/D3)^2+(G3/F3)^2)*H3"

# --- Selection change ---
$null = $ws.Range("E5").Select()
